$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.524.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.481.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.865.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.515.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.762"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.535.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("E21").Value = "  +6.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.972.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.722.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
